# Auto-generated edit script: apply cell value updates per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "308.93"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.17%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.08"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-1.33%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.236"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.06%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07664"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.63%"
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.620"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.47%"
$ws.Range("B7").Value = "BTSEToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "2.486"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.43%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9188"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.61%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1243"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "15.02%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1837"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "4.28%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09105"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.67%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04292"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.82%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1050"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.01%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001260"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.00%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005737"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.06%"
$ws.Range("B16").Value = "UpBots"
$ws.Range("C16").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.007498"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2,392.15%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.354"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.21%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.319"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1.53%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.229"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "10.18%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2894"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "7.88%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04064"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-2.52%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001265"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3.49%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004139"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.15%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-2.07%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02452"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "1.76%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05297"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "1.87%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007853"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.18%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1313"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.11%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006827"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.79%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001914"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.29%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008362"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "6.13%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.26%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006667"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.96%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.26%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2055"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "1,853.97%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-2.44%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.26%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.26%"
